# edit.ps1 - apply the "skip write fault at finishing exec due to notification
# from cli service" change to CLI_terminal_comm_protocol.docx
#
# Summary of changes (see unified diff):
#  1) Merge several runs that made up one sentence about "Окончание работы
#     команды с последующим завершением клиентской сессии..." into a single
#     run (pure text consolidation, no visible change).
#  2) Merge several runs that made up "Результатом обработки запроса на
#     выполнение команды..." into a single run.
#  3) Merge several runs that made up the long sentence about "...клиентская
#     сессия должна быть завершена (также как и выполнение терминального
#     клиента)..." into a single run.
#  4) Merge two runs "...терминального клиента. " + "Ответ " into one run.
#  5) Wrap the "Запрос на завершение работы {exit}" run pair in a bookmark
#     (__DdeLink__271_1487678931).
#  6) Append a brand-new list paragraph at the end of the document body
#     explaining that the terminal client may, in response, generate its own
#     {exit} request.
#  7) Register two new (unused) character styles ListLabel47 / ListLabel48
#     that LibreOffice mints alongside list-using paragraphs.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Окончание работы команды с последующим завершением..." run merge
# ---------------------------------------------------------------------
$t1 = "Окончание работы команды с последующим завершением клиентской сессии (и выполнения терминального клиента) "
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "Результатом обработки запроса на выполнение команды..." run merge
# ---------------------------------------------------------------------
$t2 = "Результатом обработки запроса на выполнение команды на CLI сервисе будут несколько ответов "
$d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "...клиентская сессия должна быть завершена (также как и выполнение
#    терминального клиента), то будет прислан ответ " run merge
# ---------------------------------------------------------------------
$t3 = " (несколько ответов означает, что будет прислано любое количество или 0). После завершения выполнения команды будет прислан ответ о завершении. Если после выполнения команды клиентская сессия должна быть завершена (также как и выполнение терминального клиента), то будет прислан ответ "
$d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "...терминального клиента. " + "Ответ " run merge
# ---------------------------------------------------------------------
$t4 = " означает, что пользователю не удалось войти в систему (скорее всего из-за того, что он ввел неверное имя пользователя и/или пароль), но у пользователя есть еще минимум одна возможность ввести данные аутентификации через данный экземпляр терминального клиента. Ответ "
$d.Content.Find.Execute($t4, $true, $false, $false, $false, $false, $true, 1, $false, $t4, 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Bookmark around "Запрос на завершение работы " + "{exit}"
#    (the document already has bookmarks 0..9; the new one becomes 10)
# ---------------------------------------------------------------------
$bmStartRng = $d.Content
$bmStartRng.Find.Execute("Запрос на завершение работы ", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$bmStart = $bmStartRng.Start

$bmEndRng = $d.Range($bmStartRng.End, $d.Content.End)
$bmEndRng.Find.Execute("{exit}", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$bmEnd = $bmEndRng.End

$d.Bookmarks.Add("__DdeLink__271_1487678931", $d.Range($bmStart, $bmEnd)) | Out-Null

# ---------------------------------------------------------------------
# 6) New trailing paragraph about the terminal client re-issuing {exit}
# ---------------------------------------------------------------------
$insertPoint = $d.Range($d.Content.End, $d.Content.End)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">При получении уведомления о завершении клиентской сессии терминальный клиент в ответ может сгенерировать запрос на завершение работы </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{exit}</w:t></w:r><w:r><w:rPr/><w:t>. Будет ли в действительности генерироваться этот запрос зависит от реализации терминального клиента.</w:t></w:r></w:p>'
$insertPoint.InsertXML($newParaXml)

# Recolor the trailing (non-bold) sentence explicitly as "not bold" to match
# the existing document convention used for the other non-bold runs in this
# list (b/bCs explicit-false on runs following a bold {placeholder}).
$tail = ". Будет ли в действительности генерироваться этот запрос зависит от реализации терминального клиента."
$tailFind = $d.Content.Find
$tailFind.ClearFormatting()
$tailFind.Replacement.ClearFormatting()
$tailFind.Replacement.Font.Bold = $false
$tailFind.Replacement.Font.BoldBi = $false
$tailFind.Execute($tail, $true, $false, $false, $false, $false, $true, 1, $false, $tail, 2) | Out-Null

# ---------------------------------------------------------------------
# 7) New (unused) character styles minted by LibreOffice for list labels
# ---------------------------------------------------------------------
$lbl47 = $d.Styles.Add("ListLabel 47", 2)
$lbl47.Font.NameBi = "Symbol"

$lbl48 = $d.Styles.Add("ListLabel 48", 2)
$lbl48.Font.NameBi = "OpenSymbol"

Write-Output "edit complete"
